$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at sheet row 30 (shifts DESCUENTOS..GRUPO_ASEGURADOR down to 31..36)
$ws.Rows.Item(30).Insert()

# Populate the new row with the PRIMA_FRACCIONADA column entry
$ws.Range("A30").Value = "PRIMA_FRACCIONADA"

# Grow the table (tbl_polizas) so the new row becomes part of it
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F36"))

# Update the active cell selection to match the authored edit
$ws.Range("B30").Select()
